$d = $word.ActiveDocument

# The run inside the "RunField" content control ("Something new") was
# corrupted with five duplicate <w:rPr> children, which violates the
# OOXML schema for w:r (it allows only a single rPr). Re-apply the run's
# existing (italic, red) character formatting through the Word object
# model so the run gets rewritten with a single, consolidated <w:rPr>.

$text = $d.Content.Text
$target = "Something new"
$start = $text.IndexOf($target)

if ($start -ge 0) {
    $range = $d.Range($start, $start + $target.Length)
    # Re-assert the formatting that is already present on the run; this
    # forces the engine to regenerate a clean, single <w:rPr> element
    # instead of the five stacked duplicates.
    $range.Font.Italic = 1
    $range.Font.Color = 255
}
